$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2
$ws.Range("F2").Value = 31
$ws.Range("G2").Value = "adam"
$ws.Range("I2").Value = 64
$ws.Range("J2").Value = 18.39320491063797
$ws.Range("K2").Value = 529.1887752371133
$ws.Range("L2").Value = 23.00410344345359
$ws.Range("M2").Value = 0.1155928170780067

# Add new row 3
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "DNN"
$ws.Range("D3").Value = 30
$ws.Range("E3").Value = 60
$ws.Range("F3").Value = 31
$ws.Range("G3").Value = "<keras.src.optimizers.legacy.adam.Adam object at 0x7b85ac652860>"
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 32
$ws.Range("J3").Value = 20.37185309281508
$ws.Range("K3").Value = 663.5352690490155
$ws.Range("L3").Value = 25.75917834576669
$ws.Range("M3").Value = 0.129579387327666

# Add new row 4
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "DNN"
$ws.Range("D4").Value = 40
$ws.Range("E4").Value = 60
$ws.Range("F4").Value = 31
$ws.Range("G4").Value = "adam"
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 64
$ws.Range("J4").Value = 24.52132068249605
$ws.Range("K4").Value = 974.0897322202044
$ws.Range("L4").Value = 31.21041063844249
$ws.Range("M4").Value = 0.1510193738722635

# Add new row 5
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "DNN"
$ws.Range("D5").Value = 40
$ws.Range("E5").Value = 60
$ws.Range("F5").Value = 31
$ws.Range("G5").Value = "<keras.src.optimizers.legacy.adam.Adam object at 0x7b8623c04a90>"
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 32
$ws.Range("J5").Value = 18.61304731960103
$ws.Range("K5").Value = 575.3971132186766
$ws.Range("L5").Value = 23.98743657039402
$ws.Range("M5").Value = 0.1188064917914747
